$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.851.53'
$ws.Range("E2").Value = '  -0.15%  '

$ws.Range("D3").Value = '3.523.10'
$ws.Range("E3").Value = '  -0.01%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.65%  '

$ws.Range("D7").Value = '3.520.94'
$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  +0.83%  '

$ws.Range("E10").Value = '  +1.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.13'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.35%  '

$ws.Range("E12").Value = '  +0.12%  '

$ws.Range("D13").Value = '4.124.20'
$ws.Range("E13").Value = '  -0.04%  '

$ws.Range("E14").Value = '  +0.68%  '

$ws.Range("E15").Value = '  +0.62%  '

$ws.Range("E16").Value = '  +0.01%  '

$ws.Range("D17").Value = '3.522.75'
$ws.Range("E17").Value = '  -0.20%  '

$ws.Range("D18").Value = '64.922.69'
$ws.Range("E18").Value = '  -0.04%  '

$ws.Range("E19").Value = '  -0.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '391.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.95%  '

$ws.Range("E23").Value = '  +0.74%  '

$ws.Range("D24").Value = '3.666.18'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.23'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("E27").Value = '  +0.47%  '

$ws.Range("E28").Value = '  +20.53%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.03%  '

$ws.Range("E31").Value = '  +1.54%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.64%  '

$ws.Range("D33").Value = '3.527.76'
$ws.Range("E33").Value = '  -0.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.09'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.85%  '

$ws.Range("E36").Value = '  +1.02%  '

$ws.Range("E37").Value = '  +6.16%  '

$ws.Range("E38").Value = '  +2.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '169.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.45%  '

$ws.Range("E40").Value = '  +0.42%  '

$ws.Range("E41").Value = '  +3.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.822'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.65'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.47%  '

$ws.Range("E44").Value = '  +3.32%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.42'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.56%  '

$ws.Range("E48").Value = '  -0.21%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.97%  '

$ws.Range("D50").Value = '2.404.74'
$ws.Range("E50").Value = '  +0.07%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.895'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.26%  '
